$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 42
$ws.Range("F2").Value = 0.8501760114036275
$ws.Range("G2").Value = 0.1291565590113278
$ws.Range("H2").Value = 1447907280199614000
$ws.Range("I2").Value = -0.153584975624098
$ws.Range("J2").Value = -0.007732616684352993
$ws.Range("K2").Value = 0.6150591811991899
$ws.Range("L2").Value = 1.121563450511849
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.1500191543621032

# Row 3
$ws.Range("B3").Value = 42
$ws.Range("F3").Value = 0.0076414631584662
$ws.Range("G3").Value = 0.08447820241948407
$ws.Range("H3").Value = -10612932943871520
$ws.Range("I3").Value = 0.4516945429508499
$ws.Range("J3").Value = -0.001589774102599029
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0.9895578760275494
$ws.Range("N3").Value = 0.1500191543621032

# Row 4
$ws.Range("B4").Value = 42
$ws.Range("F4").Value = 0.0003565913521299812
$ws.Range("G4").Value = 0.01223163318571148
$ws.Range("H4").Value = 47472369842610310
$ws.Range("I4").Value = 0.8722687363902777
$ws.Range("J4").Value = -0.0008500983906317937
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0.9991113085980893
$ws.Range("N4").Value = 0.1500191543621032

# Row 5
$ws.Range("B5").Value = 42
$ws.Range("F5").Value = 0.9249944634425169
$ws.Range("G5").Value = 0.1423367181911414
$ws.Range("H5").Value = -1842791083890418000
$ws.Range("I5").Value = 0.6437610425204559
$ws.Range("J5").Value = 0.02764200895731305
$ws.Range("K5").Value = 0.6635441657487544
$ws.Range("L5").Value = 1.227159110107581
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0.07580466036226756

# Row 6
$ws.Range("B6").Value = 42
$ws.Range("F6").Value = 0.00004653781548358664
$ws.Range("G6").Value = 0.003121851568369667
$ws.Range("H6").Value = 22891814706713270
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = -0.0002222715939667628
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0.9997778271495223
$ws.Range("N6").Value = 0.07580466036226756

# Row 7
$ws.Range("B7").Value = 42
$ws.Range("F7").Value = 0.01097788015030839
$ws.Range("G7").Value = 0.1062416093779839
$ws.Range("H7").Value = 7340467833711607
$ws.Range("I7").Value = -1.892465721900112
$ws.Range("J7").Value = -0.01067933391970033
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0.9862252832703844
$ws.Range("N7").Value = 0.07580466036226756
